$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83:155 down to 84:156.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record's data.
$ws.Cells.Item(83, 1).Value2 = 5
$ws.Cells.Item(83, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(83, 3).Value2 = "Maule"
$ws.Cells.Item(83, 4).Value2 = 45240
$ws.Cells.Item(83, 5).Value2 = 7
$ws.Cells.Item(83, 6).Value2 = 100112022
$ws.Cells.Item(83, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(83, 8).Value2 = "Sin especificar"
$ws.Cells.Item(83, 9).Value2 = "Primera"
$ws.Cells.Item(83, 10).Value2 = 500
$ws.Cells.Item(83, 11).Value2 = 20000
$ws.Cells.Item(83, 12).Value2 = 22000
$ws.Cells.Item(83, 13).Value2 = 21200
$ws.Cells.Item(83, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(83, 15).Value2 = "Región del Maule"
$ws.Cells.Item(83, 16).Value2 = 848
$ws.Cells.Item(83, 17).Value2 = 25
$ws.Cells.Item(83, 18).Value2 = "Hortaliza"
